# LOT2017.xlsx rebuild: "Docentes responsaveis" row (old row 13, which only
# held the "4873328 - Fernando Segato" value in B/C with no label in A) is
# removed entirely, shifting every row below it up by one. The long
# paragraph texts that used to live under "Objetivos:", "Programa resumido:",
# "Programa:", "Metodo:" and "Bibliografia:" are replaced with short values,
# and the remaining evaluation/bibliography cells shift into the row above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "4873328 - Fernando Segato" row (old row 13); this shifts
# rows 14-26 up to become rows 13-25 and keeps row heights / column-A labels
# aligned with their correct rows automatically.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) - replace the long mission paragraph with the teacher id.
$ws.Range("B10").Value = "4873328 - Fernando Segato"
$ws.Range("C10").Value = "4873328 - Fernando Segato"

# Row 13 (Programa resumido:) - replace the long summary paragraph.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:) - replace the long syllabus paragraph.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# Row 18 (Metodo:) - replace the evaluation-method text.
$ws.Range("B18").Value = "4873328 - Fernando Segato"
$ws.Range("C18").Value = "4873328 - Fernando Segato"

# Row 19 (Criterio:) - now holds what used to be the "Metodo:" text.
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas (P1 e P2)."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas (P1 e P2)."

# Row 20 (Norma de recuperação:) - now holds what used to be the "Criterio:" text.
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1x1 + P2x2)/3"

# Row 21 (Bibliografia:) - now holds what used to be the "Norma de recuperação:" text.
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
